$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$ws.Range("B5").Value = "C:\Users\Josue Mk\Documents\UiPath\RPA_Practica_Git_Github\Input\"
$ws.Range("B5").ReadingOrder = 0
